# Update the "Förändrad" (Changed) date in column C for all data rows
# (rows 2-120) from 2023-10-05 (45204) to 2023-10-08 (45207).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C120").Value = 45207
